$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45179 to 45180
# (2023-09-10 -> 2023-09-11) for every data row (rows 2 through 141).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 141 }

$ws.Range("C2:C$lastRow").Value = 45180
